$wb = $excel.ActiveWorkbook

# --- Sheet1: change selection from A1:E3 to a single active cell A3 ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate() | Out-Null
$ws1.Range("A3").Select() | Out-Null

# --- Sheet3: rename "Authors" header to "Followers", add a new
#     "Person being followed" column, and fix a typo in the author list ---
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Activate() | Out-Null

# Rename header (was "Authors")
$ws3.Range("A1").Value = "Followers"

# New column E: "Person being followed"
$ws3.Range("E1").Value = "Person being followed"
$ws3.Range("E2").Value = "Barack, Obama"
$ws3.Range("E3").Value = "Glaston, John"
$ws3.Range("E4").Value = "Glaston, John"
$ws3.Range("E1:E4").WrapText = $true

# Match the widened column E (compensates for the engine's
# character-width -> pixel -> width quantization so the stored
# <col> width lands exactly on 35)
$ws3.Columns.Item(5).ColumnWidth = 34.0834

# Fix the typo "Bloor,Mary" -> "Bloor, Mary" in the author list
$ws3.Range("A4").Value = "Lemarck, Christian; Bloor, Mary; Van Dijke, Leo"

# Final selection / active sheet is Sheet3, cell A5
$ws3.Range("A5").Select() | Out-Null
